$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.856.69"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "3.435.72"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.88"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.56"
$ws.Range("E6").Value = "  +5.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("D8").Value = "3.433.76"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.176"
$ws.Range("E10").Value = "  +5.95%  "
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.42"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "3.986.71"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.52"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "3.442.08"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "66.716.32"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.99"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.05"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.97"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.93"
$ws.Range("E24").Value = "  +10.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.17"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.44"
$ws.Range("E26").Value = "  +3.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.95"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.91"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.91"
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.47"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.97"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.96"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "582.26"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("E36").Value = "  +5.87%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  +5.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.57"
$ws.Range("E40").Value = "  +2.65%  "
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "3.127.44"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("E44").Value = "  +6.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0424"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("E46").Value = "  +19.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.65"
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("E51").Value = "  +5.21%  "
